$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CTLT")

# Row 4: Inventory
$ws.Range("B4").Value = 461000000.0
$ws.Range("C4").Value = 384000000.0
$ws.Range("D4").Value = 324000000.0
$ws.Range("E4").Value = 297000000.0
$ws.Range("F4").Value = 250000000.0

# Row 13: Accounts Payable
$ws.Range("B13").Value = 357000000.0
$ws.Range("C13").Value = 329000000.0
$ws.Range("D13").Value = 321000000.0
$ws.Range("E13").Value = 276000000.0
$ws.Range("F13").Value = 243000000.0

# Row 23: Long Term Tax Liability (Deferred)
$ws.Range("B23").Value = 22000000.0
$ws.Range("C23").Value = 27000000.0
$ws.Range("D23").Value = 45000000.0
$ws.Range("E23").Value = 51000000.0
$ws.Range("F23").Value = 41000000.0
